$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SPY Walls")

# Row 57: fill in C:H (previously empty)
$ws.Range("C57").Value = 588.1
$ws.Range("D57").Value = 595.54
$ws.Range("E57").Value = 588.1
$ws.Range("F57").Value = 594.85
$ws.Range("G57").Value = 67928536
$ws.Range("H57").Value = 0.40515468548855943

# Row 58: fill in most cells (A,B date markers + I..II data), leave C:H empty
$ws.Range("A58").Value = 45796
$ws.Range("B58").Value = 45797
$ws.Range("I58").Value = 0.18140000000000001
$ws.Range("J58").Value = 4.0999999999999996
$ws.Range("K58").Value = 600
$ws.Range("L58").Value = 319192200
$ws.Range("M58").Value = 2061
$ws.Range("N58").Value = 4378
$ws.Range("O58").Value = 6439
$ws.Range("P58").Value = 0.11626311362979302
$ws.Range("Q58").Value = 0.17496784793250422
$ws.Range("R58").Value = 45797
$ws.Range("S58").Value = 0.17496784793250422
$ws.Range("T58").Value = 45800
$ws.Range("U58").Value = 0.10646206056035498
$ws.Range("V58").Value = 45828
$ws.Range("W58").Value = 0.35509012679024732
$ws.Range("X58").Value = 12.333333333333334
$ws.Range("Y58").Value = 605
$ws.Range("Z58").Value = 175916455
$ws.Range("AA58").Value = -8531
$ws.Range("AB58").Value = 41
$ws.Range("AC58").Value = 8572
$ws.Range("AD58").Value = 0.064076110873058215
$ws.Range("AE58").Value = 0.034134332046919263
$ws.Range("AF58").Value = 45800
$ws.Range("AG58").Value = 0.062343554938752752
$ws.Range("AH58").Value = 45807
$ws.Range("AI58").Value = 0.24946293796695881
$ws.Range("AJ58").Value = 45828
$ws.Range("AK58").Value = 0.36607669055721376
$ws.Range("AL58").Value = 15.666666666666666
$ws.Range("AM58").Value = 595
$ws.Range("AN58").Value = 129298855
$ws.Range("AO58").Value = 41037
$ws.Range("AP58").Value = 44947
$ws.Range("AQ58").Value = 85984
$ws.Range("AR58").Value = 0.047096036403982085
$ws.Range("AS58").Value = 0.23575440794405481
$ws.Range("AT58").Value = 45797
$ws.Range("AU58").Value = 0.23575440794405481
$ws.Range("AV58").Value = 45798
$ws.Range("AW58").Value = 0.15662192830223554
$ws.Range("AX58").Value = 45807
$ws.Range("AY58").Value = 0.1637896509586059
$ws.Range("AZ58").Value = 4.666666666666667
$ws.Range("BA58").Value = 610
$ws.Range("BB58").Value = 99381810
$ws.Range("BC58").Value = 160
$ws.Range("BD58").Value = 892
$ws.Range("BE58").Value = 1052
$ws.Range("BF58").Value = 0.03619900069226159
$ws.Range("BG58").Value = 0
$ws.Range("BH58").Value = 45807
$ws.Range("BI58").Value = 0.13091320883722918
$ws.Range("BJ58").Value = 45814
$ws.Range("BK58").Value = 0.10899106875774539
$ws.Range("BL58").Value = 45828
$ws.Range("BM58").Value = 0.23205204905773258
$ws.Range("BN58").Value = 20.333333333333332
$ws.Range("BO58").Value = 615
$ws.Range("BP58").Value = 84944415
$ws.Range("BQ58").Value = 6256
$ws.Range("BR58").Value = 130
$ws.Range("BS58").Value = 6386
$ws.Range("BT58").Value = 0.030940299209571205
$ws.Range("BU58").Value = 0
$ws.Range("BV58").Value = 45828
$ws.Range("BW58").Value = 0.15000789732109135
$ws.Range("BX58").Value = 45856
$ws.Range("BY58").Value = 0.2627885097411739
$ws.Range("BZ58").Value = 45919
$ws.Range("CA58").Value = 0.15995852189618112
$ws.Range("CB58").Value = 71.666666666666671
$ws.Range("CC58").Value = 575
$ws.Range("CD58").Value = -71471925
$ws.Range("CE58").Value = 0.038940822115318088
$ws.Range("CF58").Value = -1207
$ws.Range("CG58").Value = 3080
$ws.Range("CH58").Value = 4287
$ws.Range("CI58").Value = 0
$ws.Range("CJ58").Value = 45807
$ws.Range("CK58").Value = 0.20116551919754519
$ws.Range("CL58").Value = 45828
$ws.Range("CM58").Value = 0.33938268739846833
$ws.Range("CN58").Value = 45856
$ws.Range("CO58").Value = 0.077058353317346118
$ws.Range("CP58").Value = 34.333333333333336
$ws.Range("CQ58").Value = 585
$ws.Range("CR58").Value = -55011645
$ws.Range("CS58").Value = 0.029972589687713432
$ws.Range("CT58").Value = 9417
$ws.Range("CU58").Value = 20515
$ws.Range("CV58").Value = 29932
$ws.Range("CW58").Value = 0.06967510502989957
$ws.Range("CX58").Value = 45800
$ws.Range("CY58").Value = 0.12676522793999992
$ws.Range("CZ58").Value = 45807
$ws.Range("DA58").Value = 0.24255301730714129
$ws.Range("DB58").Value = 45828
$ws.Range("DC58").Value = 0.10196738026219333
$ws.Range("DD58").Value = 15.666666666666666
$ws.Range("DE58").Value = 570
$ws.Range("DF58").Value = -42560760
$ws.Range("DG58").Value = 0.0231888393135171
$ws.Range("DH58").Value = -1198
$ws.Range("DI58").Value = 10172
$ws.Range("DJ58").Value = 11370
$ws.Range("DK58").Value = 0
$ws.Range("DL58").Value = 45807
$ws.Range("DM58").Value = 0.27058817242226996
$ws.Range("DN58").Value = 45828
$ws.Range("DO58").Value = 0.16952041983305011
$ws.Range("DP58").Value = 45838
$ws.Range("DQ58").Value = 0.10964504440952962
$ws.Range("DR58").Value = 28.333333333333332
$ws.Range("DS58").Value = 565
$ws.Range("DT58").Value = -40344390
$ws.Range("DU58").Value = 0.021981270468663296
$ws.Range("DV58").Value = -977
$ws.Range("DW58").Value = 3661
$ws.Range("DX58").Value = 4638
$ws.Range("DY58").Value = 0
$ws.Range("DZ58").Value = 45807
$ws.Range("EA58").Value = 0.16968284434695183
$ws.Range("EB58").Value = 45828
$ws.Range("EC58").Value = 0.16168540352801389
$ws.Range("ED58").Value = 45856
$ws.Range("EE58").Value = 0.35737135545196963
$ws.Range("EF58").Value = 34.333333333333336
$ws.Range("EG58").Value = 580
$ws.Range("EH58").Value = -28450160
$ws.Range("EI58").Value = 0.015500808460277769
$ws.Range("EJ58").Value = 887
$ws.Range("EK58").Value = 12407
$ws.Range("EL58").Value = 13294
$ws.Range("EM58").Value = 0.054795467863395314
$ws.Range("EN58").Value = 45800
$ws.Range("EO58").Value = 0.13628439048003815
$ws.Range("EP58").Value = 45807
$ws.Range("EQ58").Value = 0.31359531790834544
$ws.Range("ER58").Value = 45828
$ws.Range("ES58").Value = 0.13179665804061161
$ws.Range("ET58").Value = 15.666666666666666
$ws.Range("EU58").Value = 600
$ws.Range("EV58").Value = 393670200
$ws.Range("EW58").Value = 2061
$ws.Range("EX58").Value = 4378
$ws.Range("EY58").Value = 6439
$ws.Range("EZ58").Value = 0.085938653436710033
$ws.Range("FA58").Value = 356431200
$ws.Range("FB58").Value = 0.12982711077151471
$ws.Range("FC58").Value = 0.17496784793250422
$ws.Range("FD58").Value = 45797
$ws.Range("FE58").Value = 0.17496784793250422
$ws.Range("FF58").Value = 45800
$ws.Range("FG58").Value = 0.10646206056035498
$ws.Range("FH58").Value = 45828
$ws.Range("FI58").Value = 0.35509012679024732
$ws.Range("FJ58").Value = 12.333333333333334
$ws.Range("FK58").Value = -37239000
$ws.Range("FL58").Value = 0.020289327239364694
$ws.Range("FM58").Value = 0.016482719729316038
$ws.Range("FN58").Value = 45828
$ws.Range("FO58").Value = 0.16589059856601951
$ws.Range("FP58").Value = 45884
$ws.Range("FQ58").Value = 0.14747442197695965
$ws.Range("FR58").Value = 45919
$ws.Range("FS58").Value = 0.12677032143720293
$ws.Range("FT58").Value = 81
$ws.Range("FU58").Value = 595
$ws.Range("FV58").Value = 334234705
$ws.Range("FW58").Value = 41037
$ws.Range("FX58").Value = 44947
$ws.Range("FY58").Value = 85984
$ws.Range("FZ58").Value = 0.072963817122850591
$ws.Range("GA58").Value = 231766780
$ws.Range("GB58").Value = 0.084419128909638888
$ws.Range("GC58").Value = 0.23575440794405481
$ws.Range("GD58").Value = 45797
$ws.Range("GE58").Value = 0.23575440794405481
$ws.Range("GF58").Value = 45798
$ws.Range("GG58").Value = 0.15662192830223554
$ws.Range("GH58").Value = 45807
$ws.Range("GI58").Value = 0.1637896509586059
$ws.Range("GJ58").Value = 4.666666666666667
$ws.Range("GK58").Value = -102467925
$ws.Range("GL58").Value = 0.05582870812491416
$ws.Range("GM58").Value = 0.31152338646459365
$ws.Range("GN58").Value = 45797
$ws.Range("GO58").Value = 0.31152338646459365
$ws.Range("GP58").Value = 45800
$ws.Range("GQ58").Value = 0.087983044450251138
$ws.Range("GR58").Value = 45828
$ws.Range("GS58").Value = 0.2866010510118166
$ws.Range("GT58").Value = 12.333333333333334
$ws.Range("GU58").Value = 590
$ws.Range("GV58").Value = 329627100
$ws.Range("GW58").Value = -670
$ws.Range("GX58").Value = 18984
$ws.Range("GY58").Value = 19654
$ws.Range("GZ58").Value = 0.07195797169876654
$ws.Range("HA58").Value = 151835320
$ws.Range("HB58").Value = 0.055304757015290425
$ws.Range("HC58").Value = 0.13623575858370768
$ws.Range("HD58").Value = 45797
$ws.Range("HE58").Value = 0.13623575858370768
$ws.Range("HF58").Value = 45807
$ws.Range("HG58").Value = 0.19247089544119247
$ws.Range("HH58").Value = 45828
$ws.Range("HI58").Value = 0.14328846542425044
$ws.Range("HJ58").Value = 14.666666666666666
$ws.Range("HK58").Value = -177791780
$ws.Range("HL58").Value = 0.096868218934158673
$ws.Range("HM58").Value = 0.31570773406959535
$ws.Range("HN58").Value = 45797
$ws.Range("HO58").Value = 0.31570773406959535
$ws.Range("HP58").Value = 45798
$ws.Range("HQ58").Value = 0.075014435425529796
$ws.Range("HR58").Value = 45800
$ws.Range("HS58").Value = 0.23437157780860285
$ws.Range("HT58").Value = 2.3333333333333335
$ws.Range("HU58").Value = 595
$ws.Range("HV58").Value = 141801
$ws.Range("HW58").Value = 212194
$ws.Range("HX58").Value = 2745429655.5
$ws.Range("HY58").Value = -1835398461.5
$ws.Range("HZ58").Value = 910031194
$ws.Range("IA58").Value = 1.4958221405809977
$ws.Range("IB58").Value = 4580828117
$ws.Range("IC58").Value = 0.17057982335118468
$ws.Range("ID58").Value = 45797
$ws.Range("IE58").Value = 0.17057982335118468
$ws.Range("IF58").Value = 45807
$ws.Range("IG58").Value = 0.11823665550557919
$ws.Range("IH58").Value = 45828
$ws.Range("II58").Value = 0.19190208376028442

# Update the selected/active cell in the bottom-right pane to C60
$ws.Range("C60").Select()